$wb = $excel.ActiveWorkbook

# --- "Daily" sheet: update row 2 (clear/cloudy sky GHI forecast values) ---
$daily = $wb.Worksheets.Item("Daily")
$daily.Range("G2").Value = 2531.71
$daily.Range("H2").Value = 5706.21
$daily.Range("J2").Value = 681.29
$daily.Range("L2").Value = 681.29

# --- "Hourly" sheet: update forecasted hourly values ---
$hourly = $wb.Worksheets.Item("Hourly")

# row 9
$hourly.Range("K9").Value = 0.31
$hourly.Range("M9").Value = 0.31

# row 10
$hourly.Range("I10").Value = 343.89
$hourly.Range("K10").Value = 21.14
$hourly.Range("M10").Value = 21.14

# row 11
$hourly.Range("H11").Value = 207.48
$hourly.Range("I11").Value = 591.51
$hourly.Range("J11").Value = 66.91

# row 12
$hourly.Range("H12").Value = 325.68
$hourly.Range("I12").Value = 702.5700000000001

# row 13
$hourly.Range("I13").Value = 756.23
$hourly.Range("K13").Value = 101.41
$hourly.Range("M13").Value = 101.41

# row 14
$hourly.Range("K14").Value = 115.5
$hourly.Range("M14").Value = 115.5

# row 15
$hourly.Range("K15").Value = 119.96
$hourly.Range("M15").Value = 119.96

# row 16
$hourly.Range("I16").Value = 714.48
$hourly.Range("K16").Value = 97.11
$hourly.Range("M16").Value = 97.11

# row 17
$hourly.Range("H17").Value = 228.95
$hourly.Range("I17").Value = 615.65
$hourly.Range("K17").Value = 61.01
$hourly.Range("M17").Value = 61.01

# row 18
$hourly.Range("I18").Value = 402.73
$hourly.Range("K18").Value = 29.62
$hourly.Range("M18").Value = 29.62

# row 19
$hourly.Range("K19").Value = 1.93
$hourly.Range("M19").Value = 1.93
